{"js": "// Update the multiplication problems in the practice table to a new\n// generated set (output regenerated at c986bee).\nconst replacements = [\n  [\"12\u00d714=\", \"96\u00d764=\"],\n  [\"81\u00d714=\", \"85\u00d741=\"],\n  [\"89\u00d779=\", \"17\u00d754=\"],\n  [\"20\u00d799=\", \"45\u00d718=\"],\n  [\"84\u00d774=\", \"17\u00d759=\"],\n  [\"67\u00d729=\", \"88\u00d733=\"],\n  [\"63\u00d717=\", \"30\u00d742=\"],\n  [\"25\u00d761=\", \"77\u00d729=\"],\n  [\"11\u00d779=\", \"71\u00d762=\"],\n  [\"44\u00d748=\", \"30\u00d765=\"],\n  [\"89\u00d737=\", \"96\u00d790=\"],\n  [\"40\u00d726=\", \"23\u00d739=\"],\n  [\"20\u00d757=\", \"58\u00d740=\"],\n  [\"42\u00d742=\", \"14\u00d731=\"],\n  [\"42\u00d736=\", \"87\u00d796=\"],\n  [\"74\u00d724=\", \"83\u00d719=\"],\n  [\"13\u00d797=\", \"18\u00d775=\"],\n  [\"47\u00d795=\", \"12\u00d782=\"],\n  [\"74\u00d790=\", \"79\u00d775=\"],\n  [\"20\u00d755=\", \"44\u00d721=\"],\n  [\"87\u00d735=\", \"26\u00d768=\"],\n  [\"31\u00d780=\", \"88\u00d784=\"],\n  [\"16\u00d785=\", \"66\u00d722=\"],\n  [\"93\u00d712=\", \"91\u00d735=\"],\n  [\"76\u00d775=\", \"74\u00d745=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the multiplication problems in the practice table to a new\n# generated set (output regenerated at c986bee).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"12\u00d714=\"; New = \"96\u00d764=\" },\n    @{ Old = \"81\u00d714=\"; New = \"85\u00d741=\" },\n    @{ Old = \"89\u00d779=\"; New = \"17\u00d754=\" },\n    @{ Old = \"20\u00d799=\"; New = \"45\u00d718=\" },\n    @{ Old = \"84\u00d774=\"; New = \"17\u00d759=\" },\n    @{ Old = \"67\u00d729=\"; New = \"88\u00d733=\" },\n    @{ Old = \"63\u00d717=\"; New = \"30\u00d742=\" },\n    @{ Old = \"25\u00d761=\"; New = \"77\u00d729=\" },\n    @{ Old = \"11\u00d779=\"; New = \"71\u00d762=\" },\n    @{ Old = \"44\u00d748=\"; New = \"30\u00d765=\" },\n    @{ Old = \"89\u00d737=\"; New = \"96\u00d790=\" },\n    @{ Old = \"40\u00d726=\"; New = \"23\u00d739=\" },\n    @{ Old = \"20\u00d757=\"; New = \"58\u00d740=\" },\n    @{ Old = \"42\u00d742=\"; New = \"14\u00d731=\" },\n    @{ Old = \"42\u00d736=\"; New = \"87\u00d796=\" },\n    @{ Old = \"74\u00d724=\"; New = \"83\u00d719=\" },\n    @{ Old = \"13\u00d797=\"; New = \"18\u00d775=\" },\n    @{ Old = \"47\u00d795=\"; New = \"12\u00d782=\" },\n    @{ Old = \"74\u00d790=\"; New = \"79\u00d775=\" },\n    @{ Old = \"20\u00d755=\"; New = \"44\u00d721=\" },\n    @{ Old = \"87\u00d735=\"; New = \"26\u00d768=\" },\n    @{ Old = \"31\u00d780=\"; New = \"88\u00d784=\" },\n    @{ Old = \"16\u00d785=\"; New = \"66\u00d722=\" },\n    @{ Old = \"93\u00d712=\"; New = \"91\u00d735=\" },\n    @{ Old = \"76\u00d775=\"; New = \"74\u00d745=\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($pair.Old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
